$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("transferability")
$add = $wb.Worksheets.Item("additivity")

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Sheet1"

$src.Range("B18:G43").Copy()
$newSheet.Range("A1").PasteSpecial()

$newSheet.Move($add)
